$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReadMe")
$ws.Range("A1").Value = "Hello"
